# Apply updated ESFUERZO (column C) values to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2 = 0.18616
    3 = 0.18599
    4 = 0.20461
    5 = 0.22108
    6 = 0.22889
    7 = 0.20063
    8 = 0.23031
    9 = 0.21167
    10 = 0.22398
    11 = 0.23698
    12 = 0.21439
    13 = 0.2627
    14 = 0.219
    15 = 0.24329
    16 = 0.20393
    17 = 0.20787
    18 = 0.20283
    19 = 0.18659
    20 = 0.20847
    21 = 0.26376
    22 = 0.27376
    23 = 0.2089
    24 = 0.2043
    25 = 0.20828
    26 = 0.20405
    27 = 0.24447
    28 = 0.27593
    29 = 0.34555
    30 = 0.3127
    31 = 0.27338
    32 = 0.26192
    33 = 0.28826
    34 = 0.22127
    35 = 0.21426
    36 = 0.26167
    37 = 0.31332
    38 = 0.28524
    39 = 0.24511
    40 = 0.20692
    41 = 0.21628
    42 = 0.22478
    43 = 0.26344
    44 = 0.23459
    45 = 0.3267
    46 = 0.28926
    47 = 0.26061
    48 = 0.27042
    49 = 0.22961
    50 = 0.1884
    51 = 0.17757
    52 = 0.2331
    53 = 0.23172
    54 = 0.24164
    55 = 0.19157
    56 = 0.22639
    57 = 0.17455
    58 = 0.18802
    59 = 0.21972
    60 = 0.20033
    61 = 0.26
    62 = 0.21013
    63 = 0.25068
    64 = 0.21477
    65 = 0.22451
    66 = 0.25012
    67 = 0.24327
    68 = 0.24589
    69 = 0.30311
    70 = 0.2919
    71 = 0.28324
    72 = 0.3042
    73 = 0.25898
    74 = 0.19685
    75 = 0.20264
    76 = 0.19472
    77 = 0.2491
    78 = 0.25858
    79 = 0.24201
    80 = 0.22153
    81 = 0.22214
    82 = 0.14885
    83 = 0.15075
    84 = 0.1707
    85 = 0.20759
    86 = 0.13324
    87 = 0.14612
    88 = 0.23042
    89 = 0.17023
    90 = 0.19521
    91 = 0.20742
    92 = 0.21076
    93 = 0.22225
    94 = 0.22499
    95 = 0.1963
    96 = 0.19727
    97 = 0.16787
    98 = 0.24554
    99 = 0.23881
    100 = 0.24886
    101 = 0.3092
    102 = 0.30767
    103 = 0.28096
    104 = 0.28408
    105 = 0.25704
    106 = 0.18735
    107 = 0.19969
    108 = 0.20276
    109 = 0.2516
    110 = 0.2161
    111 = 0.1863
    112 = 0.22096
    113 = 0.17663
    114 = 0.20913
    115 = 0.19679
    116 = 0.18693
    117 = 0.23062
    118 = 0.22051
    119 = 0.20971
    120 = 0.21746
    121 = 0.23508
    122 = 0.26486
    123 = 0.25916
    124 = 0.28922
    125 = 0.30675
    126 = 0.3
    127 = 0.27447
    128 = 0.28897
    129 = 0.30333
    130 = 0.1849
    131 = 0.19406
    132 = 0.19752
    133 = 0.22071
    134 = 0.21285
    135 = 0.16038
    136 = 0.17007
    137 = 0.18371
    138 = 0.12473
    139 = 0.27134
    140 = 0.14648
    141 = 0.29473
    142 = 0.10773
    143 = 0.22918
    144 = 0.14573
    145 = 0.23118
    146 = 0.15788
    147 = 0.10065
    148 = 0.14439
    149 = 0.09970000000000001
    150 = 0.1573
    151 = 0.19403
    152 = 0.18761
    153 = 0.20066
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}
